$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(4, 9).Value = 'sd'
$ws.Cells.Item(4, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(10, 9).Value = 'ba'
$ws.Cells.Item(10, 10).Value = 'Appreciation'
$ws.Cells.Item(14, 9).Value = 'sv'
$ws.Cells.Item(14, 10).Value = 'Statement-opinion'
$ws.Cells.Item(16, 9).Value = 'sd'
$ws.Cells.Item(16, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(17, 9).Value = 'b'
$ws.Cells.Item(17, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(23, 9).Value = 'b'
$ws.Cells.Item(23, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(26, 9).Value = 'sv'
$ws.Cells.Item(26, 10).Value = 'Statement-opinion'
$ws.Cells.Item(41, 9).Value = 'b'
$ws.Cells.Item(41, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(43, 9).Value = 'sv'
$ws.Cells.Item(43, 10).Value = 'Statement-opinion'
$ws.Cells.Item(49, 9).Value = 'aa'
$ws.Cells.Item(49, 10).Value = 'Agree/Accept'
$ws.Cells.Item(55, 9).Value = 'sv'
$ws.Cells.Item(55, 10).Value = 'Statement-opinion'
$ws.Cells.Item(59, 9).Value = 'aa'
$ws.Cells.Item(59, 10).Value = 'Agree/Accept'
$ws.Cells.Item(65, 9).Value = 'sd'
$ws.Cells.Item(65, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(71, 9).Value = 'sd'
$ws.Cells.Item(71, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(82, 9).Value = 'sd'
$ws.Cells.Item(82, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(88, 9).Value = 'sd'
$ws.Cells.Item(88, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(106, 9).Value = 'sv'
$ws.Cells.Item(106, 10).Value = 'Statement-opinion'
$ws.Cells.Item(115, 9).Value = 'sv'
$ws.Cells.Item(115, 10).Value = 'Statement-opinion'
$ws.Cells.Item(130, 9).Value = 'sd'
$ws.Cells.Item(130, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(142, 9).Value = 'sd'
$ws.Cells.Item(142, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(149, 9).Value = 'sd'
$ws.Cells.Item(149, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(152, 9).Value = 'aa'
$ws.Cells.Item(152, 10).Value = 'Agree/Accept'
$ws.Cells.Item(166, 9).Value = 'ba'
$ws.Cells.Item(166, 10).Value = 'Appreciation'
$ws.Cells.Item(171, 9).Value = 'sv'
$ws.Cells.Item(171, 10).Value = 'Statement-opinion'
$ws.Cells.Item(174, 9).Value = 'sd'
$ws.Cells.Item(174, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(186, 9).Value = 'b'
$ws.Cells.Item(186, 10).Value = 'Acknowledge (Backchannel)'
$ws.Cells.Item(187, 9).Value = 'sv'
$ws.Cells.Item(187, 10).Value = 'Statement-opinion'
$ws.Cells.Item(200, 9).Value = 'sd'
$ws.Cells.Item(200, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(206, 9).Value = 'sv'
$ws.Cells.Item(206, 10).Value = 'Statement-opinion'
$ws.Cells.Item(213, 9).Value = 'sv'
$ws.Cells.Item(213, 10).Value = 'Statement-opinion'
$ws.Cells.Item(219, 9).Value = 'sd'
$ws.Cells.Item(219, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(224, 9).Value = 'sv'
$ws.Cells.Item(224, 10).Value = 'Statement-opinion'
$ws.Cells.Item(227, 9).Value = 'sd'
$ws.Cells.Item(227, 10).Value = 'Statement-non-opinion'
$ws.Cells.Item(232, 9).Value = 'sv'
$ws.Cells.Item(232, 10).Value = 'Statement-opinion'
$ws.Cells.Item(233, 9).Value = 'ba'
$ws.Cells.Item(233, 10).Value = 'Appreciation'
$ws.Cells.Item(237, 9).Value = '%'
$ws.Cells.Item(237, 10).Value = 'Uninterpretable'
